$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.827.12"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "2.660.83"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'601.49"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "'154.96"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.546"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").Value = "2.658.52"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("E10").Value = "  +11.46%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "'27.81"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("E15").Value = "  +5.06%  "
$ws.Range("D16").Value = "3.142.20"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "68.710.26"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "2.656.86"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").Value = "'11.46"
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("D20").Value = "'367.52"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "'7.50"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").Value = "'73.27"
$ws.Range("E25").Value = "  +8.65%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'9.91"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.767.58"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "'0.0000105"
$ws.Range("E29").Value = "  +3.19%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'581.70"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").Value = "'7.97"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("D37").Value = "'1.53"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("D38").Value = "'160.57"
$ws.Range("E38").Value = "  +2.65%  "
$ws.Range("D39").Value = "'19.29"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "'0.367"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "'5.36"
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").Value = "'2.66"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").Value = "'17.63"
$ws.Range("E44").Value = "  +4.83%  "
$ws.Range("D45").Value = "0.0₆0322"
$ws.Range("E45").Value = "  +8.43%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'40.65"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "'156.37"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "'1.71"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'21.95"
$ws.Range("E51").Value = "  -0.32%  "
